$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the first occurrence of $oldSub inside the text of the
# paragraph at $paraIndex with $newSub, by operating on a precise character
# Range (keeps the rest of the paragraph's runs/formatting untouched).
# ---------------------------------------------------------------------------
function Replace-InParagraphText([int]$paraIndex, [string]$oldSub, [string]$newSub) {
    $p = $d.Paragraphs.Item($paraIndex).Range
    $text = $p.Text
    $pos = $text.IndexOf($oldSub)
    if ($pos -lt 0) {
        return $false
    }
    $absStart = $p.Start + $pos
    $absEnd = $absStart + $oldSub.Length
    $sub = $d.Range($absStart, $absEnd)
    $sub.Text = $newSub
    return $true
}

# ---------------------------------------------------------------------------
# 1) Sprint Completion Date:  07/25/2016  ->  07/24/2016
# ---------------------------------------------------------------------------
Replace-InParagraphText 7 "25" "24" | Out-Null

# ---------------------------------------------------------------------------
# 2) Revision Number:  1  ->  3
# ---------------------------------------------------------------------------
Replace-InParagraphText 8 "1" "3" | Out-Null

# ---------------------------------------------------------------------------
# 3) Revision Date:  07/17/2016  ->  07/18/2016
# ---------------------------------------------------------------------------
Replace-InParagraphText 9 "17" "18" | Out-Null

# ---------------------------------------------------------------------------
# 4) Re-purpose the "Tasks Completed This Sprint" Heading1 paragraph (now
#    paragraph 10) into the first of a new set of section headings, then
#    insert the rest of the new headings (with a blank paragraph between
#    each), and finally thin the big block of blank "spacing" paragraphs
#    that used to sit below it down to just the two that still belong to
#    the document (immediately above "Images of prototype created").
# ---------------------------------------------------------------------------
$tasksHeading = $d.Paragraphs.Item(10)
$tasksHeading.Range.Text = "actions to stop doing"

# Remove 7 of the 9 blank "spacing" paragraphs that used to pad out this
# section - 2 of them stay, right before the next Heading1 ("Images of
# prototype created").
for ($k = 0; $k -lt 7; $k++) {
    $blank = $d.Paragraphs.Item(11)
    $blank.Range.Delete()
}

function Insert-BareParagraphAfter([int]$afterIndex) {
    $p = $d.Paragraphs.Item($afterIndex)
    $p.Range.InsertParagraphAfter()
    $newIndex = $afterIndex + 1
    $np = $d.Paragraphs.Item($newIndex)
    $np.Range.Style = "Normal"
    return $newIndex
}

function Insert-Heading1After([int]$afterIndex, [string]$text) {
    $p = $d.Paragraphs.Item($afterIndex)
    $p.Range.InsertParagraphAfter()
    $newIndex = $afterIndex + 1
    $np = $d.Paragraphs.Item($newIndex)
    $np.Range.Style = "Heading1"
    $np.Range.Text = $text
    return $newIndex
}

$idx = 10
$idx = Insert-BareParagraphAfter $idx
$idx = Insert-Heading1After $idx "Actions to start doing"
$idx = Insert-BareParagraphAfter $idx
$idx = Insert-Heading1After $idx "Actions to Keep doing"
$idx = Insert-BareParagraphAfter $idx
$idx = Insert-Heading1After $idx "Tasks Completed This Sprint"
$idx = Insert-BareParagraphAfter $idx
$idx = Insert-Heading1After $idx "Tasks not completed this Sprint"
$idx = Insert-BareParagraphAfter $idx
$idx = Insert-Heading1After $idx "Work Completion Rate "

Write-Output "Done."
